# Fill out the data dictionary: add Source/SourceHeading entries for the
# M-M bridge tables, add an Explanation column (I) describing the bridge
# tables, and rename the old Explanation header (H1) to "Restriction" now
# that Explanation has moved to column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Cells.Item(1, 8).Value  = "Restriction"
$ws.Cells.Item(1, 9).Value  = "Explanation"

# --- Character_Actor bridge (rows 6-7): explanation ------------------
$ws.Cells.Item(6, 9).Value  = "M-M bridge from character to actor"
$ws.Cells.Item(7, 9).Value  = "M-M bridge from character to actor"

# --- Actor table (rows 8-11): Source / SourceHeading ------------------
$ws.Cells.Item(8,  6).Value = "actors"
$ws.Cells.Item(8,  7).Value = "gender"
$ws.Cells.Item(9,  6).Value = "actors"
$ws.Cells.Item(9,  7).Value = "fname"
$ws.Cells.Item(10, 6).Value = "actors"
$ws.Cells.Item(10, 7).Value = "lname"
$ws.Cells.Item(11, 6).Value = "actors"
$ws.Cells.Item(11, 7).Value = "idactors"

# --- Actor_Aka_Names table (rows 12-14) -------------------------------
$ws.Cells.Item(12, 6).Value = "aka_names"
$ws.Cells.Item(12, 7).Value = "idactors"
$ws.Cells.Item(13, 6).Value = "aka_names"
$ws.Cells.Item(13, 7).Value = "idaka_names"
$ws.Cells.Item(14, 6).Value = "aka_names"
$ws.Cells.Item(14, 7).Value = "name"

# --- Character_Movie bridge (rows 15-16): explanation -----------------
$ws.Cells.Item(15, 9).Value = "M-M bridge from character to movie"
$ws.Cells.Item(16, 9).Value = "M-M bridge from character to movie"

# --- Movie table (rows 17-19) -----------------------------------------
$ws.Cells.Item(17, 6).Value = "moves"
$ws.Cells.Item(17, 7).Value = "idmovies"
$ws.Cells.Item(18, 6).Value = "movies"
$ws.Cells.Item(18, 7).Value = "title"
$ws.Cells.Item(19, 6).Value = "movies"
$ws.Cells.Item(19, 7).Value = "year"

# --- Movie_Keyword table (rows 20-21) ---------------------------------
$ws.Cells.Item(20, 6).Value = "movies_keywords"
$ws.Cells.Item(20, 7).Value = "idmovies"
$ws.Cells.Item(21, 6).Value = "movies_keywords"
$ws.Cells.Item(21, 7).Value = "idkeywords"

# --- Keyword table (rows 22-23) ---------------------------------------
$ws.Cells.Item(22, 6).Value = "keywords"
$ws.Cells.Item(22, 7).Value = "idkeywords"
$ws.Cells.Item(23, 6).Value = "keywords"
$ws.Cells.Item(23, 7).Value = "keyword"

# --- Aka_Movie table (rows 24-27) -------------------------------------
$ws.Cells.Item(24, 6).Value = "aka_titles"
$ws.Cells.Item(24, 7).Value = "idmovies"
$ws.Cells.Item(25, 6).Value = "aka_titles"
$ws.Cells.Item(25, 7).Value = "idaka_titles"
$ws.Cells.Item(26, 6).Value = "aka_titles"
$ws.Cells.Item(26, 7).Value = "year"
$ws.Cells.Item(27, 6).Value = "aka_titles"
$ws.Cells.Item(27, 7).Value = "title"

# --- Character_Series bridge (rows 28-29): explanation ----------------
$ws.Cells.Item(28, 9).Value = "M-M bridge from character to series"
$ws.Cells.Item(29, 9).Value = "M-M bridge from character to series"

# --- Series table (rows 30-34) ----------------------------------------
$ws.Cells.Item(30, 6).Value = "series"
$ws.Cells.Item(30, 7).Value = "name"
$ws.Cells.Item(31, 6).Value = "series"
$ws.Cells.Item(31, 7).Value = "idmovies"
$ws.Cells.Item(32, 6).Value = "series"
$ws.Cells.Item(32, 7).Value = "season"
$ws.Cells.Item(33, 6).Value = "series"
$ws.Cells.Item(33, 7).Value = "number"
$ws.Cells.Item(34, 6).Value = "series"
$ws.Cells.Item(34, 7).Value = "idseries"

# --- Genre table (rows 35-36) ------------------------------------------
$ws.Cells.Item(35, 6).Value = "genres"
$ws.Cells.Item(35, 7).Value = "idgenres"
$ws.Cells.Item(36, 6).Value = "genres"
$ws.Cells.Item(36, 7).Value = "genre"

# --- Column I width (best match reachable in this engine's rounding) --
$ws.Columns.Item(9).ColumnWidth = 9.8

# --- Selection / view: author ended editing at I30, scrolled to top ---
$ws.Range("I30").Select()
